$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 has no gaps (A:O all populated) -- a single block copy preserves styles without padding.
$ws.Range("A2:O2").Copy($ws.Range("A4:O4"))
$ws.Range("A2:O2").Copy($ws.Range("A6:O6"))

# Row 3 has gaps (I3, N3 are absent) -- copy cell-by-cell so the destination
# doesn't get padded with spurious empty cells in those columns.
$row3Cols = @("A","B","C","D","E","F","G","H","J","K","L","M","O")
foreach ($col in $row3Cols) {
    $ws.Range("${col}3").Copy($ws.Range("${col}5"))
    $ws.Range("${col}3").Copy($ws.Range("${col}7"))
}

# Update selection to match final state
$ws.Range("A6:XFD7").Select()
